$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Grab references to the existing sheets
# ----------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("RQ1")
$ws2 = $wb.Worksheets.Item("RQ2")
$ws3 = $wb.Worksheets.Item("RQ3")
$ws4 = $wb.Worksheets.Item("RQ4")

# ----------------------------------------------------------------------------
# RQ1: update rows 5-8 with the new (reduced revision count) figures and
# insert a new row 12 holding the previous Jabref totals, annotated with a
# note about the revision count used.
# ----------------------------------------------------------------------------

# Row 5 (Freecol)
$ws1.Range("C5").Value = 14
$ws1.Range("D5").Value = 3
$ws1.Range("E5").Formula = "=C5/D5"
$ws1.Range("F5").Value = 18
$ws1.Range("G5").Value = 3
$ws1.Range("H5").Formula = "=F5/G5"

# Row 6 (Carol)
$ws1.Range("C6").Value = 0
$ws1.Range("D6").Value = 2
$ws1.Range("E6").Formula = "=C6/D6"
$ws1.Range("F6").Value = 8
$ws1.Range("G6").Value = 2
$ws1.Range("H6").Formula = "=F6/G6"

# Row 7 (Jabref) - replaced with the reduced-revision figures
$ws1.Range("C7").Value = 0
$ws1.Range("D7").Value = 2
$ws1.Range("F7").Value = 6
$ws1.Range("G7").Value = 2

# Row 8 totals recalculate automatically through the existing SUM formulas

# New row 12: carries the previous (full revision count) Jabref figures
$ws1.Range("A7").Copy()
$ws1.Range("A12").PasteSpecial(-4122)
$ws1.Range("A12").Value = "Jabref"
$ws1.Range("B12").Value = 2798
$ws1.Range("C12").Value = 117
$ws1.Range("D12").Value = 65
$ws1.Range("E12").Formula = "=C12/D12"
$ws1.Range("F12").Value = 316
$ws1.Range("G12").Value = 65
$ws1.Range("H12").Formula = "=F12/G12"

# ----------------------------------------------------------------------------
# RQ2: new Freecol / Carol / Jabref figures
# ----------------------------------------------------------------------------
$ws2.Range("B5").Value = 14
$ws2.Range("C5").Value = 14
$ws2.Range("E5").Value = 18
$ws2.Range("F5").Value = 548

$ws2.Range("B6").Value = 0
$ws2.Range("C6").Value = 0
$ws2.Range("E6").Value = 8
$ws2.Range("F6").Value = 9

$ws2.Range("B7").Value = 0
$ws2.Range("C7").Value = 0
$ws2.Range("E7").Value = 6
$ws2.Range("F7").Value = 147

# ----------------------------------------------------------------------------
# RQ3: new Freecol / Carol / Jabref figures
# ----------------------------------------------------------------------------
$ws3.Range("B5").Value = 3
$ws3.Range("C5").Value = 3
$ws3.Range("E5").Value = 3
$ws3.Range("F5").Value = 4

$ws3.Range("B6").Value = 1
$ws3.Range("C6").Value = 1
$ws3.Range("E6").Value = 2
$ws3.Range("F6").Value = 2

$ws3.Range("B7").Value = 1
$ws3.Range("C7").Value = 1
$ws3.Range("E7").Value = 2
$ws3.Range("F7").Value = 4

# ----------------------------------------------------------------------------
# RQ4: new Freecol / Carol / Jabref figures
# ----------------------------------------------------------------------------
$ws4.Range("B5").Value = 3
$ws4.Range("C5").Value = 1
$ws4.Range("E5").Value = 3
$ws4.Range("F5").Value = 1

$ws4.Range("B6").Value = 1
$ws4.Range("C6").Value = 1
$ws4.Range("E6").Value = 2
$ws4.Range("F6").Value = 1

$ws4.Range("B7").Value = 1
$ws4.Range("C7").Value = 1
$ws4.Range("E7").Value = 2
$ws4.Range("F7").Value = 1

# ----------------------------------------------------------------------------
# Add the new RQ5 sheet (LOCs-based results) after RQ4
# ----------------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws5 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws5.Name = "RQ5"

# Reuse the header / label formatting from RQ2 (same visual layout)
$ws2.Range("A1:G2").Copy()
$ws5.Range("A1").PasteSpecial(-4122)
$ws2.Range("A3:A8").Copy()
$ws5.Range("A3").PasteSpecial(-4122)

$ws5.Range("B1").Value = "Regular Clones"
$ws5.Range("E1").Value = "Micro Clones"

$ws5.Range("B2").Value = "Rep CFs LOCs"
$ws5.Range("C2").Value = "CFs LOCs"
$ws5.Range("D2").Value = "%"
$ws5.Range("E2").Value = "Rep CFs LOCs"
$ws5.Range("F2").Value = "CFs LOCs"
$ws5.Range("G2").Value = "%"

$ws5.Range("A3").Value = "Ctags"
$ws5.Range("A4").Value = "Brlcad"
$ws5.Range("A5").Value = "Freecol"
$ws5.Range("A6").Value = "Carol"
$ws5.Range("A7").Value = "Jabref"
$ws5.Range("A8").Value = "Total"

$ws5.Range("B3").Value = 437894
$ws5.Range("C3").Value = 5177512
$ws5.Range("D3").Formula = "=B3/C3*100"
$ws5.Range("E3").Value = 2925129
$ws5.Range("F3").Value = 17036371
$ws5.Range("G3").Formula = "=E3/F3*100"

$ws5.Range("B4").Value = 1363688
$ws5.Range("C4").Value = 12733290
$ws5.Range("D4").Formula = "=B4/C4*100"
$ws5.Range("E4").Value = 2317473
$ws5.Range("F4").Value = 19098646
$ws5.Range("G4").Formula = "=E4/F4*100"

$ws5.Range("B5").Value = 202006
$ws5.Range("C5").Value = 3030064
$ws5.Range("D5").Formula = "=B5/C5*100"
$ws5.Range("E5").Value = 563836
$ws5.Range("F5").Value = 6576597
$ws5.Range("G5").Formula = "=E5/F5*100"

$ws5.Range("B6").Value = 0
$ws5.Range("C6").Value = 17116
$ws5.Range("D6").Formula = "=B6/C6*100"
$ws5.Range("E6").Value = 9148
$ws5.Range("F6").Value = 47747
$ws5.Range("G6").Formula = "=E6/F6*100"

$ws5.Range("B7").Value = 0
$ws5.Range("C7").Value = 480701
$ws5.Range("D7").Formula = "=B7/C7*100"
$ws5.Range("E7").Value = 46320
$ws5.Range("F7").Value = 705773
$ws5.Range("G7").Formula = "=E7/F7*100"

$ws5.Range("B8").Formula = "=SUM(B3:B7)"
$ws5.Range("C8").Formula = "=SUM(C3:C7)"
$ws5.Range("D8").Formula = "=B8/C8*100"
$ws5.Range("E8").Formula = "=SUM(E3:E7)"
$ws5.Range("F8").Formula = "=SUM(F3:F7)"
$ws5.Range("G8").Formula = "=E8/F8*100"

$ws5.Range("B1:D1").Merge()
$ws5.Range("E1:G1").Merge()

$ws5.Columns.Item(2).ColumnWidth = 12.6
$ws5.Columns.Item(4).ColumnWidth = 8.31
$ws5.Columns.Item(5).ColumnWidth = 12.74

# ----------------------------------------------------------------------------
# RQ1: widen column I so the new annotation text fits, and annotate row 12
# ----------------------------------------------------------------------------
$ws1.Columns.Item(9).ColumnWidth = 19.7
$ws1.Range("I12").Value = "For 1545 revisions"

# ----------------------------------------------------------------------------
# View / selection state: restore per-sheet selections and make RQ5 the
# active tab, matching the author's final screen state.
# ----------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("H5").Select()

$ws2.Activate()
$ws2.Range("G8").Select()

$ws3.Activate()
$ws3.Range("F5").Select()

$ws4.Activate()
$ws4.Range("F5").Select()

$ws5.Activate()
$ws5.Range("G5").Select()
